# sn: update the pretas forms
#
# Renames the "school" fields (d_cluster_name / d_cluster_id) labels to
# "village" labels, bumps the form_id/form_title to the "_v2" / "V2"
# revision, and re-selects the sheets/cells the author left active.

$wb = $excel.ActiveWorkbook

$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- content changes -------------------------------------------------

# Row 5 = d_cluster_name: "4. Entrer le nom de l'école" -> "4. Nom du village"
$survey.Range("C5").Value = "4. Nom du village"

# Row 6 = d_cluster_id: "5. Entrer le code de l'école" -> "5. Code du village"
$survey.Range("C6").Value = "5. Code du village"

# Row 9 (begin repeat) internal name: sn_lf_f_2407 -> sn_lf_f_2407_2
$survey.Range("B9").Value = "sn_lf_f_2407_2"

# settings sheet: form_title / form_id bumped to the V2 / _v2 revision
$settings.Range("A2").Value = "(2024 Juillet) 2. Pre-TAS - Formulaire Résultat FTS V2"
$settings.Range("B2").Value = "sn_lf_pretas_20407_2_fts_v2"

# --- view / selection changes -----------------------------------------

# survey sheet is no longer the tab left selected; its frozen pane scrolls
# back to the top and the cursor rests on B9 instead of I30
$survey.Activate()
$survey.Range("B9").Select()

# settings becomes the active (tab-selected) sheet, cursor on B2
$settings.Activate()
$settings.Range("B2").Select()
